$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number, report week range) ---
$ws.Range("A8").Value = "Volume 29   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/7/2022  Through  11/13/2022"

# --- Crime-stat table updates ---
# Row 14
$ws.Range("N14").Value = -37.5
# Row 15
$c = $ws.Range("C15")
$c.Value = 1
$c.NumberFormat = "#,##0"
$ws.Range("L15").Value = 5.555555555555
$ws.Range("M15").Value = 58.333333333333
$ws.Range("N15").Value = -32.142857142857
# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 63.636363636363
$ws.Range("I16").Value = 182
$ws.Range("J16").Value = 102
$ws.Range("K16").Value = 78.431372549019
$ws.Range("L16").Value = 56.896551724137
$ws.Range("M16").Value = -26.016260162601
$ws.Range("N16").Value = -82.330097087378
# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = 13.793103448275
$ws.Range("I17").Value = 330
$ws.Range("J17").Value = 238
$ws.Range("K17").Value = 38.655462184873
$ws.Range("L17").Value = 50
$ws.Range("M17").Value = 64.179104477611
$ws.Range("N17").Value = 2.167182662538
# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -9.090909090909
$ws.Range("I18").Value = 135
$ws.Range("J18").Value = 93
$ws.Range("K18").Value = 45.16129032258
$ws.Range("L18").Value = 3.053435114503
$ws.Range("M18").Value = -55.882352941176
$ws.Range("N18").Value = -89.277204130262
# Row 19
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = -3.225806451612
$ws.Range("I19").Value = 329
$ws.Range("J19").Value = 229
$ws.Range("K19").Value = 43.668122270742
$ws.Range("L19").Value = 33.739837398374
$ws.Range("M19").Value = 6.818181818181
$ws.Range("N19").Value = -34.98023715415
# Row 20
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 16.666666666666
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 227
$ws.Range("J20").Value = 202
$ws.Range("K20").Value = 12.376237623762
$ws.Range("L20").Value = 69.402985074626
$ws.Range("M20").Value = -11.328125
$ws.Range("N20").Value = -92.518127883981
# Row 21
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 20.833333333333
$ws.Range("F21").Value = 116
$ws.Range("G21").Value = 98
$ws.Range("H21").Value = 18.367346938775
$ws.Range("I21").Value = 1227
$ws.Range("J21").Value = 889
$ws.Range("K21").Value = 38.020247469066
$ws.Range("L21").Value = 41.03448275862
$ws.Range("M21").Value = -8.020989505247
$ws.Range("N21").Value = -80.171299288946
# Row 24
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -28.571428571428
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 111
$ws.Range("H24").Value = 9.009009009009
$ws.Range("I24").Value = 1430
$ws.Range("J24").Value = 890
$ws.Range("K24").Value = 60.67415730337
$ws.Range("L24").Value = 72.08182912154
$ws.Range("M24").Value = 103.994293865906
# Row 25
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 11
$ws.Range("F25").Value = 49
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 32.432432432432
$ws.Range("I25").Value = 484
$ws.Range("J25").Value = 446
$ws.Range("K25").Value = 8.520179372197
$ws.Range("L25").Value = 24.742268041237
$ws.Range("M25").Value = -12.949640287769
# Row 26
$c = $ws.Range("C26")
$c.Value = 1
$c.NumberFormat = "#,##0"
$ws.Range("L26").Value = 6.666666666666
# Row 27
$c = $ws.Range("C27")
$c.Value = 3
$c.NumberFormat = "#,##0"
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 40
$ws.Range("I27").Value = 47
$ws.Range("J27").Value = 63
$ws.Range("K27").Value = -25.396825396825
$ws.Range("L27").Value = 4.444444444444
# Row 28
$ws.Range("N28").Value = -28.571428571428
# Row 29
$ws.Range("N29").Value = -56.521739130434
